$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 ---
# C7:F7 currently carry style s="5" (green fill + thick-left/right border) and value 2.
# Target: no explicit style (plain/default), value 5.
# Strip the style by pasting the (unstyled) format from C4:F4, which already has no "s".
$ws.Range("C4:F4").Copy()
$ws.Range("C7:F7").PasteSpecial(-4122)
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 5

# G7/H7 keep their existing style (s="2"), they just gain a value.
$ws.Range("G7").Value = 5
$ws.Range("H7").Value = 5

# New cells I7, J7 (plain/no style, value 5) and K7 (plain, text "Хочу 4",
# the same shared string already used by K4/K11).
$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 5
$ws.Range("K7").Value = "Хочу 4"

# --- Row 22 ---
# C22/D22/E22 keep their style (s="5"); value changes from 2 to 5.
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = 5

# G22/H22 keep their existing style (s="2"), they just gain a value.
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 5

# --- Row 25 ---
# D25 keeps its style (s="5"); value changes from 2 to 5.
$ws.Range("D25").Value = 5

# G25 keeps its existing style (s="2"), it just gains a value.
$ws.Range("G25").Value = 5

# --- New green-filled cells I22, J22, I25 (thick left border) and J25 (no border) ---
# Build each distinct look once on a "prototype" cell, then fan the format out with
# copy/paste-special so only ONE new xf entry is created per distinct look (instead of
# one per property-assignment per cell).
$protoBordered = $ws.Range("I22")
$protoBordered.Interior.Color = 5296274
$protoBordered.Borders.Item(7).LineStyle = 1
$protoBordered.Borders.Item(7).Weight = 4
$protoBordered.Borders.Item(7).Color = 0
$protoBordered.HorizontalAlignment = -4108
$protoBordered.VerticalAlignment = -4108
$protoBordered.WrapText = $true

$protoPlain = $ws.Range("J25")
$protoPlain.Interior.Color = 5296274
$protoPlain.HorizontalAlignment = -4108
$protoPlain.VerticalAlignment = -4108
$protoPlain.WrapText = $true

$ws.Range("I22").Copy()
$ws.Range("J22").PasteSpecial(-4122)
$ws.Range("I25").PasteSpecial(-4122)

# Now write the values into all four new cells.
$ws.Range("I22").Value = 5
$ws.Range("J22").Value = 5
$ws.Range("I25").Value = 5
$ws.Range("J25").Value = 5

# --- View state: frozen pane scrolled back to the top, selection moved to K9 ---
$ws.Range("K9").Select()
